# Applies the crypto price/volume refresh for the Mon Dec 18 20:12:39 UTC 2023 GitHub Actions update.
# Row 30/31 also swap (Toncoin now ranks above Monero).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.725.86"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "2.176.02"
$ws.Range("E3").Value = "  -2.82%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'238.04"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'0.612"
$ws.Range("E6").Value = "  -2.38%  "
$ws.Range("D7").Value = "'72.42"
$ws.Range("E7").Value = "  -2.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -3.39%  "
$ws.Range("D10").Value = "'40.02"
$ws.Range("E10").Value = "  -5.13%  "
$ws.Range("D11").Value = "'0.0908"
$ws.Range("E11").Value = "  -5.02%  "
$ws.Range("D12").Value = "'54.53"
$ws.Range("E12").Value = "  -3.60%  "
$ws.Range("D13").Value = "'0.100"
$ws.Range("E13").Value = "  -3.12%  "
$ws.Range("D14").Value = "'6.69"
$ws.Range("E14").Value = "  -3.61%  "
$ws.Range("D15").Value = "2.506.08"
$ws.Range("E15").Value = "  -2.63%  "
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "2.173.58"
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "'0.778"
$ws.Range("E18").Value = "  -7.28%  "
$ws.Range("D19").Value = "41.601.71"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -2.05%  "
$ws.Range("D21").Value = "'70.13"
$ws.Range("E21").Value = "  -3.60%  "
$ws.Range("D22").Value = "'5.80"
$ws.Range("E22").Value = "  -6.75%  "
$ws.Range("D23").Value = "'10.02"
$ws.Range("E23").Value = "  -11.03%  "
$ws.Range("D24").Value = "'226.03"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'10.70"
$ws.Range("E27").Value = "  -5.87%  "
$ws.Range("E28").Value = "  -10.15%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -3.75%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.18"
$ws.Range("E30").Value = "  -0.93%  "
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'170.91"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").Value = "'32.57"
$ws.Range("E33").Value = "  +9.18%  "
$ws.Range("D34").Value = "'0.0774"
$ws.Range("E34").Value = "  -3.97%  "
$ws.Range("E35").Value = "  -6.21%  "
$ws.Range("E36").Value = "  -3.63%  "
$ws.Range("E37").Value = "  -0.86%  "
$ws.Range("D38").Value = "'0.103"
$ws.Range("E38").Value = "  -6.63%  "
$ws.Range("D39").Value = "'0.0308"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "'12.01"
$ws.Range("E40").Value = "  -8.82%  "
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("E42").Value = "  -5.96%  "
$ws.Range("D43").Value = "'59.21"
$ws.Range("E43").Value = "  -8.13%  "
$ws.Range("E44").Value = "  -5.02%  "
$ws.Range("D45").Value = "'8.43"
$ws.Range("E45").Value = "  -3.43%  "
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("D47").Value = "'97.22"
$ws.Range("E47").Value = "  -7.05%  "
$ws.Range("E48").Value = "  -5.06%  "
$ws.Range("E49").Value = "  -5.15%  "
$ws.Range("E50").Value = "  -5.70%  "
$ws.Range("E51").Value = "  -2.03%  "
